$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 75985
$ws.Range("E2").Value = 212183321

$ws.Range("C3").Value = 185263
$ws.Range("D3").Value = 46891
$ws.Range("E3").Value = 647680153

$ws.Range("C4").Value = 72561
$ws.Range("D4").Value = 18394
$ws.Range("E4").Value = 358762865

$ws.Range("C5").Value = 25875
$ws.Range("E5").Value = 185361484

$ws.Range("C6").Value = 11952
$ws.Range("E6").Value = 161831790

$ws.Range("C7").Value = 3482
$ws.Range("E7").Value = 110367500

$ws.Range("C14").Value = 81776
$ws.Range("D14").Value = 20600
$ws.Range("E14").Value = 159220979

$ws.Range("C15").Value = 18882
$ws.Range("E15").Value = 52006060

$ws.Range("C16").Value = 51533
$ws.Range("D16").Value = 13435
$ws.Range("E16").Value = 177302802

$ws.Range("C17").Value = 18370
$ws.Range("D17").Value = 4884
$ws.Range("E17").Value = 85434665

$ws.Range("C18").Value = 5844
$ws.Range("E18").Value = 36592400

$ws.Range("C19").Value = 2497
$ws.Range("E19").Value = 29746817

$ws.Range("C23").Value = 18296
$ws.Range("E23").Value = 33419633

$ws.Range("C24").Value = 26910
$ws.Range("E24").Value = 84827014

$ws.Range("C25").Value = 64209
$ws.Range("D25").Value = 16324
$ws.Range("E25").Value = 238252587

$ws.Range("C26").Value = 24015
$ws.Range("E26").Value = 118306531

$ws.Range("C27").Value = 7838
$ws.Range("E27").Value = 51678377

$ws.Range("C28").Value = 3298
$ws.Range("E28").Value = 39294749

$ws.Range("C29").Value = 931
$ws.Range("E29").Value = 31438425

$ws.Range("C33").Value = 20055
$ws.Range("E33").Value = 37094194

$ws.Range("C34").Value = 14744
$ws.Range("D34").Value = 3742
$ws.Range("E34").Value = 40790377

$ws.Range("C35").Value = 43291
$ws.Range("D35").Value = 11547
$ws.Range("E35").Value = 143598482

$ws.Range("C36").Value = 16616
$ws.Range("E36").Value = 75076774

$ws.Range("C37").Value = 4942
$ws.Range("E37").Value = 29873691

$ws.Range("C42").Value = 7570
$ws.Range("E42").Value = 25940511

$ws.Range("C43").Value = 17025
$ws.Range("E43").Value = 61281443

$ws.Range("C44").Value = 7663
$ws.Range("E44").Value = 38752944

$ws.Range("C49").Value = 4776
$ws.Range("E49").Value = 8432105

$ws.Range("C50").Value = 34992
$ws.Range("D50").Value = 9297
$ws.Range("E50").Value = 99571646

$ws.Range("C51").Value = 104547
$ws.Range("D51").Value = 28383
$ws.Range("E51").Value = 362084448

$ws.Range("C52").Value = 42400
$ws.Range("E52").Value = 198093959

$ws.Range("C53").Value = 15263
$ws.Range("E53").Value = 96448232

$ws.Range("C54").Value = 6719
$ws.Range("D54").Value = 1730
$ws.Range("E54").Value = 81261639

$ws.Range("C59").Value = 38355
$ws.Range("E59").Value = 89119726

$ws.Range("C60").Value = 3066
$ws.Range("E60").Value = 5556866

$ws.Range("C61").Value = 10097
$ws.Range("E61").Value = 19706272

$ws.Range("C62").Value = 3373
$ws.Range("E62").Value = 7715311

$ws.Range("C67").Value = 11177
$ws.Range("D67").Value = 2599
$ws.Range("E67").Value = 17543939

$ws.Range("C69").Value = 4871
$ws.Range("E69").Value = 11205603

$ws.Range("C71").Value = 622
$ws.Range("E71").Value = 1667796

$ws.Range("C72").Value = 241
$ws.Range("E72").Value = 794484

$ws.Range("C75").Value = 29455
$ws.Range("E75").Value = 74760299

$ws.Range("C76").Value = 87190
$ws.Range("E76").Value = 281649850

$ws.Range("C77").Value = 34071
$ws.Range("D77").Value = 8900
$ws.Range("E77").Value = 154260377

$ws.Range("C78").Value = 11555
$ws.Range("E78").Value = 67226128

$ws.Range("C79").Value = 4821
$ws.Range("E79").Value = 55412767

$ws.Range("C80").Value = 1757
$ws.Range("E80").Value = 50867173

$ws.Range("C86").Value = 26659
$ws.Range("D86").Value = 7155
$ws.Range("E86").Value = 47685786

$ws.Range("C87").Value = 109947
$ws.Range("D87").Value = 24356
$ws.Range("E87").Value = 293254379

$ws.Range("C88").Value = 297234
$ws.Range("D88").Value = 69864
$ws.Range("E88").Value = 932198740

$ws.Range("C89").Value = 143748
$ws.Range("D89").Value = 33592
$ws.Range("E89").Value = 673921243

$ws.Range("C90").Value = 58856
$ws.Range("D90").Value = 13380
$ws.Range("E90").Value = 400903856

$ws.Range("C91").Value = 27494
$ws.Range("E91").Value = 352729299

$ws.Range("C92").Value = 7751
$ws.Range("D92").Value = 2319
$ws.Range("E92").Value = 257018468

$ws.Range("C99").Value = 7
$ws.Range("D99").Value = 6
$ws.Range("E99").Value = 1210000

$ws.Range("C100").Value = 101477
$ws.Range("D100").Value = 22512
$ws.Range("E100").Value = 185101511

$ws.Range("C101").Value = 6406
$ws.Range("E101").Value = 11288817

$ws.Range("C102").Value = 14833
$ws.Range("E102").Value = 27059813

$ws.Range("C104").Value = 1714
$ws.Range("E104").Value = 4787129

$ws.Range("C108").Value = 7140
$ws.Range("E108").Value = 10128205

$ws.Range("C110").Value = 8261
$ws.Range("E110").Value = 19059921

$ws.Range("C113").Value = 419
$ws.Range("E113").Value = 2517216

$ws.Range("C116").Value = 5762
$ws.Range("E116").Value = 8769194

$ws.Range("C119").Value = 303
$ws.Range("E119").Value = 852681

$ws.Range("C123").Value = 21310
$ws.Range("E123").Value = 61208772

$ws.Range("C124").Value = 56913
$ws.Range("E124").Value = 195251975

$ws.Range("C125").Value = 21670
$ws.Range("E125").Value = 101517187

$ws.Range("C127").Value = 3011
$ws.Range("E127").Value = 36181887

$ws.Range("C128").Value = 950
$ws.Range("E128").Value = 27448906

$ws.Range("C132").Value = 16671
$ws.Range("E132").Value = 29954229

$ws.Range("C133").Value = 60422
$ws.Range("E133").Value = 178676833

$ws.Range("C134").Value = 126818
$ws.Range("D134").Value = 34482
$ws.Range("E134").Value = 425303104

$ws.Range("C135").Value = 46130
$ws.Range("D135").Value = 12717
$ws.Range("E135").Value = 212300036

$ws.Range("C136").Value = 16174
$ws.Range("E136").Value = 101856832

$ws.Range("C137").Value = 6791
$ws.Range("E137").Value = 82217962

$ws.Range("C138").Value = 2168
$ws.Range("D138").Value = 684
$ws.Range("E138").Value = 68383745

$ws.Range("C142").Value = 6
$ws.Range("E142").Value = 613000

$ws.Range("C143").Value = 44330
$ws.Range("D143").Value = 12868
$ws.Range("E143").Value = 79469039

$ws.Range("C144").Value = 74408
$ws.Range("D144").Value = 17102
$ws.Range("E144").Value = 219497216

$ws.Range("C145").Value = 152231
$ws.Range("D145").Value = 38730
$ws.Range("E145").Value = 500394254

$ws.Range("C146").Value = 52260
$ws.Range("D146").Value = 13658
$ws.Range("E146").Value = 241748668

$ws.Range("C147").Value = 17054
$ws.Range("E147").Value = 105605833

$ws.Range("C148").Value = 7012
$ws.Range("E148").Value = 81991663

$ws.Range("C154").Value = 55888
$ws.Range("E154").Value = 95148644

$ws.Range("C155").Value = 25947
$ws.Range("E155").Value = 75909607

$ws.Range("C156").Value = 63029
$ws.Range("E156").Value = 222777510

$ws.Range("C157").Value = 25463
$ws.Range("E157").Value = 120996120

$ws.Range("C158").Value = 8029
$ws.Range("E158").Value = 51208788

$ws.Range("C159").Value = 3436
$ws.Range("E159").Value = 43778492

$ws.Range("C163").Value = 18146
$ws.Range("E163").Value = 32532200

$ws.Range("C164").Value = 73007
$ws.Range("D164").Value = 16631
$ws.Range("E164").Value = 212226005

$ws.Range("C165").Value = 169099
$ws.Range("E165").Value = 567920094

$ws.Range("C166").Value = 57573
$ws.Range("E166").Value = 279411375

$ws.Range("C167").Value = 19202
$ws.Range("E167").Value = 133580055

$ws.Range("C168").Value = 8621
$ws.Range("E168").Value = 110157218

$ws.Range("C169").Value = 2560
$ws.Range("D169").Value = 765
$ws.Range("E169").Value = 87166795

$ws.Range("C175").Value = 53209
$ws.Range("E175").Value = 94536823
